$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-28 from 45465 to 45466
$ws.Range("C2:C28").Value = 45466
